# Daily attendance processing - 2026-01-18 14:04:24
# Reorders the names/emails listed in the "Recorded By" column (G) for rows
# whose value exactly matches one of the known unordered lists, so that the
# canonical ordering used by the attendance system is restored.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row  # xlUp = -4162

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -eq "System, dnasr281@gmail.com") {
        $cell.Value = "dnasr281@gmail.com, System"
    }
    elseif ($val -eq "backup@backdoor.com, System, system") {
        $cell.Value = "system, backup@backdoor.com, System"
    }
}
